$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "30.297.55";    E = "  -0.29%  " },
    @{ Row = 3;  D = "1.859.87";     E = "  -0.85%  " },
    @{ Row = 4;  D = "0.9994";       E = "  -0.11%  " },
    @{ Row = 5;  D = "233.05";       E = "  -2.29%  " },
    @{ Row = 6;  D = "0.9993";       E = "  -0.12%  " },
    @{ Row = 7;  D = "0.4756";       E = "  -1.04%  " },
    @{ Row = 8;  D = "0.2762";       E = "  -1.85%  " },
    @{ Row = 9;  D = "0.06448";      E = "  -0.94%  " },
    @{ Row = 10; D = "1.864.09";     E = "  -0.58%  " },
    @{ Row = 11; D = "0.07433";      E = "  -0.71%  " },
    @{ Row = 12; D = "16.11";        E = "  -2.65%  " },
    @{ Row = 13; D = "5.008";        E = "  -1.13%  " },
    @{ Row = 14; D = "85.77";        E = "  -2.92%  " },
    @{ Row = 15; D = "0.6356";       E = "  -3.74%  " },
    @{ Row = 16; D = "30.270.11";    E = "  -0.25%  " },
    @{ Row = 17; D = "0.9996";       E = "  -0.09%  " },
    @{ Row = 18; D = "12.83";        E = "  -3.33%  " },
    @{ Row = 19; D = "228.45";       E = "  +3.99%  " },
    @{ Row = 20; D = "0.000007380";  E = "  -2.65%  " },
    @{ Row = 21; D = "2.097.98";     E = "  -0.82%  " },
    @{ Row = 22; D = "0.9998";       E = "  -0.09%  " },
    @{ Row = 23; D = "5.127";        E = "  -3.23%  " },
    @{ Row = 24; D = "6.053";        E = "  -1.92%  " },
    @{ Row = 25; D = "9.289";        E = "  -0.45%  " },
    @{ Row = 26; D = "167.74";       E = "  +0.01%  " },
    @{ Row = 27; D = "17.93";        E = "  -2.60%  " },
    @{ Row = 28; D = "1.866";        E = "  -4.78%  " },
    @{ Row = 29; D = "0.1023";       E = "  +9.13%  " },
    @{ Row = 30; D = "1.383";        E = "  -5.42%  " },
    @{ Row = 31; D = "4.246";        E = "  -1.31%  " },
    @{ Row = 32; D = "3.921";        E = "  -2.41%  " },
    @{ Row = 33; D = "0.04899";      E = "  -2.38%  " },
    @{ Row = 34; D = "1.153";        E = "  -4.22%  " },
    @{ Row = 35; D = "0.7325";       E = "  -1.26%  " },
    @{ Row = 36; D = "0.9994";       E = "  +0.14%  " },
    @{ Row = 37; D = "2.691";        E = "  -0.64%  " },
    @{ Row = 38; D = "0.01979";      E = "  +8.67%  " },
    @{ Row = 39; D = "2.632";        E = "  +0.64%  " },
    @{ Row = 40; D = "0.9073";       E = "  +0.43%  " },
    @{ Row = 41; D = "2.000";        E = "  -2.92%  " },
    @{ Row = 42; D = "106.07";       E = "  -0.30%  " },
    @{ Row = 43; D = "0.9950";       E = "  -0.98%  " },
    @{ Row = 44; D = "0.4125";       E = "  -3.27%  " },
    @{ Row = 45; D = "5.583";        E = "  -4.62%  " },
    @{ Row = 46; D = "7.091";        E = "  -4.24%  " },
    @{ Row = 47; D = "61.37";        E = "  -4.60%  " },
    @{ Row = 48; D = "0.1210";       E = "  -4.72%  " },
    @{ Row = 49; D = "8.829";        E = "  -0.68%  " },
    @{ Row = 50; D = "1.407";        E = "  -4.46%  " },
    @{ Row = 51; D = "33.08";        E = "  -1.72%  " }
)

foreach ($u in $updates) {
    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.Style = "Normal"
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
